# February 2022 Measles Update
# Feb 2022 update for data from Jan 2022

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 38: add new monthly summary entry for Jan 2022 (K/L/M/N columns) ---
$ws.Cells.Item(38, 11).Value = 2022          # K38

$ws.Cells.Item(38, 12).Value = "January"     # L38
$ws.Cells.Item(38, 12).HorizontalAlignment = -4131   # xlLeft

$ws.Cells.Item(38, 13).Value = "NA"          # M38
$ws.Cells.Item(38, 13).HorizontalAlignment = -4131   # xlLeft

$ws.Cells.Item(38, 14).Value = "NA"          # N38
$ws.Cells.Item(38, 14).HorizontalAlignment = -4131   # xlLeft

# --- Rows 42-53: replace month-name text in column B with real dates ---
# formatted as "mmm-yy" (numFmtId 17), and normalize column C text to
# "Not reported" (already the same text, kept for clarity/explicitness).
$monthDates = @{
    42 = 44105   # Oct 2020
    43 = 44136   # Nov 2020
    44 = 44166   # Dec 2020
    45 = 44197   # Jan 2021
    46 = 44228   # Feb 2021
    47 = 44256   # Mar 2021
    48 = 44287   # Apr 2021
    49 = 44317   # May 2021
    50 = 44348   # June 2021
    51 = 44378   # July 2021
    52 = 44409   # Aug 2021
    53 = 44440   # Sep 2021
}

foreach ($r in $monthDates.Keys) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $monthDates[$r]
    $cell.NumberFormat = "mmm-yy"
    $ws.Cells.Item($r, 3).Value = "Not reported"
}

# --- New row 58: Jan 2022 reporting period, not yet reported ---
$ws.Cells.Item(58, 1).Value = 44348
$ws.Cells.Item(58, 1).NumberFormat = "d-mmm-yy"

$ws.Cells.Item(58, 2).Value = 44562
$ws.Cells.Item(58, 2).NumberFormat = "mmm-yy"

$ws.Cells.Item(58, 3).Value = "Not reported"

# --- Update the view state to reflect where the user was working ---
$null = $ws.Range("H42").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
